# Add a new "債務" (Debt) worksheet as the last sheet in the workbook.
# We copy the existing "汽車" sheet (3rd sheet) as a template, because it
# already has the same 14-column / 2-row (header+data) layout and the
# exact cell styles (s="1" for header row & col A, s="2" for the rest)
# that the new sheet needs - this avoids creating any new style entries
# in xl/styles.xml, which must stay byte-identical to the original.

$wb = $excel.ActiveWorkbook

$srcSheet = $wb.Worksheets.Item(3)
$srcSheet.Copy([System.Reflection.Missing]::Value, $wb.Worksheets.Item($wb.Worksheets.Count))
$ws = $wb.Worksheets.Item($wb.Worksheets.Count)
$ws.Name = "債務"

# The template sheet has 2 data rows; the debt sheet only needs 1.
$ws.Rows(3).Delete()

# ---- Header row (row 1) ----
$ws.Range("B1").Value = "species"
$ws.Range("C1").Value = "debtor"
# D1 "owner" already correct (inherited from template).
$ws.Range("E1").Value = "total"
$ws.Range("F1").Value = "register_date"
$ws.Range("G1").Value = "register_reason"
# H1 "property_category", I1 "category", J1 "date", K1 "legislator_name",
# L1 "legislator_id", M1 "source_file", N1 "index" already correct.

# ---- Data row (row 2) ----
$ws.Range("A2").Value = 119
$ws.Range("B2").Value = "房屋貸款"
$ws.Range("C2").Value = "王怡心"
$ws.Range("D2").Value = "國泰世華臺北市内湖區内湖路"
$ws.Range("E2").Value = 12313577
$ws.Range("F2").Value = "96年06月20日"
$ws.Range("G2").Value = "買房子"
$ws.Range("H2").Value = "debt"
# I2 "normal", J2 "2011-11-22" already correct (inherited from template).
$ws.Range("K2").Value = "費鴻泰"
# L2 1365 already correct.
# M2 "tmp1afe1" already correct.
$ws.Range("N2").Value = 119

# Sanity check (read-only) of the cells we intentionally left untouched,
# so it is obvious they still hold the values the new sheet needs.
Write-Output "H1=" $ws.Range("H1").Text "I1=" $ws.Range("I1").Text "J1=" $ws.Range("J1").Text
Write-Output "K1=" $ws.Range("K1").Text "L1=" $ws.Range("L1").Text "M1=" $ws.Range("M1").Text "N1=" $ws.Range("N1").Text
Write-Output "D1=" $ws.Range("D1").Text
Write-Output "I2=" $ws.Range("I2").Text "J2=" $ws.Range("J2").Text "L2=" $ws.Range("L2").Text "M2=" $ws.Range("M2").Text

# Restore the originally active sheet (tab 1), same as in the source file.
$wb.Worksheets.Item(1).Activate()

Write-Output "done"
